# Commit: "Using internal column widths in pptx writer tables (#9392)"
#
# The pptx table writer now honours explicit column widths instead of
# always dividing the available width evenly across every column. For
# this deck the only visible shape that contains a table is the
# "Content Placeholder 5" graphicFrame on slide 1 - its two grid
# columns grow from 2501900 EMU (197pt) to 2514600 EMU (198pt) each.
#
# Apply that using the real Table object model: locate the table,
# then set each column's width explicitly (PowerPoint measures
# Column.Width in points; 2514600 EMU == 198pt).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTable) {
        $tbl = $shape.Table
        for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
            $tbl.Columns.Item($c).Width = 198
        }
    }
}
